# Chris found that the date wasn't updated between different beta gal
# experiments. Fix the "date" column (column K) in the per-well summary
# tables on Sheet2 and Sheet3 so it reflects the correct run date
# (2018-07-27) instead of the stale 2018-07-24 value that was left over
# from a previous workbook.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Sheet2", "Sheet3")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($r = 2; $r -le 29; $r++) {
        $ws.Range("K$r").Value = "2018-07-27"
    }
}
